# Updated symbol list on Tue Jan  3 10:58:25 UTC 2023 with GitHub Actions
#
# Price (column D) and Volume(1h) (column E) cells hold numeric-looking
# text (e.g. "245.83", "-0.45%") that must stay plain text, exactly like
# the original workbook. Prefixing the value with a leading apostrophe
# forces Excel/COM to store it as text instead of silently coercing it
# to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 : BNB ---
$ws.Range("D2").Value = "'245.83"
$ws.Range("E2").Value = "'-0.45%"

# --- Row 3 : OKB ---
$ws.Range("D3").Value = "'30.22"
$ws.Range("E3").Value = "'0.39%"

# --- Row 4 : HuobiToken ---
$ws.Range("D4").Value = "'5.157"
$ws.Range("E4").Value = "'-0.33%"

# --- Row 5 : Cronos ---
$ws.Range("D5").Value = "'0.05763"
$ws.Range("E5").Value = "'0.53%"

# --- Row 6 : KuCoinToken ---
$ws.Range("D6").Value = "'6.670"
$ws.Range("E6").Value = "'1.09%"

# --- Row 7 : GateToken ---
$ws.Range("D7").Value = "'3.260"
$ws.Range("E7").Value = "'6.24%"

# --- Row 8 : MXToken ---
$ws.Range("D8").Value = "'0.8499"
$ws.Range("E8").Value = "'-1.16%"

# --- Row 9 : FTXToken ---
$ws.Range("D9").Value = "'0.8572"
$ws.Range("E9").Value = "'-2.71%"

# --- Row 10 : WazirX ---
$ws.Range("D10").Value = "'0.1385"
$ws.Range("E10").Value = "'1.26%"

# --- Row 11 : MandalaExchangeToken ---
$ws.Range("D11").Value = "'0.07080"
$ws.Range("E11").Value = "'0.13%"

# --- Row 12 : BitrueCoin ---
$ws.Range("D12").Value = "'0.03262"
$ws.Range("E12").Value = "'13.86%"

# --- Row 13 : BitMartToken ---
$ws.Range("D13").Value = "'0.09368"
$ws.Range("E13").Value = "'-0.26%"

# --- Row 14 : BitForexToken ---
$ws.Range("D14").Value = "'0.001525"
$ws.Range("E14").Value = "'0.33%"

# --- Row 15 ---
$ws.Range("D15").Value = "'0.0005948"
$ws.Range("E15").Value = "'-94.23%"

# --- Row 16 ---
$ws.Range("D16").Value = "'0.005925"
$ws.Range("E16").Value = "'-0.99%"

# --- Row 17 ---
$ws.Range("D17").Value = "'3.529"
$ws.Range("E17").Value = "'0.86%"

# --- Row 18 ---
$ws.Range("D18").Value = "'2.217"
$ws.Range("E18").Value = "'-2.15%"

# --- Row 19 ---
$ws.Range("D19").Value = "'0.3124"
$ws.Range("E19").Value = "'-1.87%"

# --- Row 20 ---
$ws.Range("D20").Value = "'0.03407"
$ws.Range("E20").Value = "'4.15%"

# --- Row 21 ---
$ws.Range("E21").Value = "'0.60%"

# --- Row 22 ---
$ws.Range("D22").Value = "'3.506"
$ws.Range("E22").Value = "'-0.66%"

# --- Row 23 (was ZBToken, now CoinExToken) ---
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04129"
$ws.Range("E23").Value = "'-0.29%"

# --- Row 24 (was CoinExToken, now ZBToken) ---
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").Value = "'0.1409"
$ws.Range("E24").Value = "'2.17%"

# --- Row 25 ---
$ws.Range("E25").Value = "'0.87%"

# --- Row 26 ---
$ws.Range("E26").Value = "'-7.80%"

# --- Row 27 ---
$ws.Range("E27").Value = "'-0.81%"

# --- Row 40 : IDEX ---
$ws.Range("D40").Value = "'0.03755"
$ws.Range("E40").Value = "'-0.96%"

# --- Row 41 : BKEXToken ---
$ws.Range("E41").Value = "'0.09%"

# --- Row 42 : CEJI ---
$ws.Range("D42").Value = "'0.002200"
$ws.Range("E42").Value = "'0.03%"

# --- Row 43 : KickToken ---
$ws.Range("D43").Value = "'0.002949"
$ws.Range("E43").Value = "'-47.31%"

# --- Row 44 : LocalTraders ---
$ws.Range("D44").Value = "'0.008933"
$ws.Range("E44").Value = "'-11.00%"

# --- Row 45 : CoinLion ---
$ws.Range("D45").Value = "'0.00005477"
$ws.Range("E45").Value = "'7.72%"

# --- Row 46 : Kangarootoken ---
$ws.Range("E46").Value = "'0.03%"

# --- Row 47 : CoinbaseStockToken ---
$ws.Range("E47").Value = "'-20.22%"

# --- Row 49 : CryptobidCoin ---
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.02%"

# --- Row 50 : SpecialPowerGold ---
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.03%"
